$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (shifts existing rows 17-23 down to 18-24)
$ws.Rows.Item(17).Insert()

# Populate the new row with the "ar" (Académie Royale) translation entry
$ws.Range("A17").Value = "ar"
$ws.Range("B17").Value = "Kunstsammlung der Académie Royale de Peinture et de Sculpture"
$ws.Range("C17").Value = "La collection d'art de l’Académie royale de peinture et de sculpture"
$ws.Range("D17").Value = "The art collection of the Académie Royale de Peinture et de Sculpture"
$ws.Rows.Item(17).RowHeight = 23.6

# Mirror the author's final selection (the newly added row)
$ws.Range("A17:D17").Select() | Out-Null
